$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 107
$ws.Cells.Item(107, 8).Value = 909.5417
$ws.Cells.Item(107, 9).Value = 884.9474
$ws.Cells.Item(107, 11).Value = 884.9474
$ws.Cells.Item(107, 13).Value = 1035.0526

# Row 138
$ws.Cells.Item(138, 8).Value = 3511826.8
$ws.Cells.Item(138, 9).Value = 8598.5
$ws.Cells.Item(138, 10).Value = 3776221.5
$ws.Cells.Item(138, 11).Value = 25795.5
$ws.Cells.Item(138, 12).Value = 11328664.5
$ws.Cells.Item(138, 13).Value = -20655.5
$ws.Cells.Item(138, 14).Value = -11338944.5

# Row 141
$ws.Cells.Item(141, 8).Value = 2816.5588
$ws.Cells.Item(141, 9).Value = 2458.7666
$ws.Cells.Item(141, 10).Value = 5500
$ws.Cells.Item(141, 11).Value = 7376.2998
$ws.Cells.Item(141, 12).Value = 16500
$ws.Cells.Item(141, 13).Value = -2196.2998
$ws.Cells.Item(141, 14).Value = -26860


$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 9834.55
$ws.Cells.Item(32, 9).Value = 5599.763
$ws.Cells.Item(32, 10).Value = 23244.709
$ws.Cells.Item(32, 11).Value = 5599.763
$ws.Cells.Item(32, 12).Value = 23244.709
$ws.Cells.Item(32, 13).Value = -5312.763
$ws.Cells.Item(32, 14).Value = -23818.709


$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 296675.16
$ws.Cells.Item(31, 9).Value = 38473.742
$ws.Cells.Item(31, 10).Value = 794635
$ws.Cells.Item(31, 11).Value = 38473.742
$ws.Cells.Item(31, 12).Value = 794635
$ws.Cells.Item(31, 13).Value = -38178.742
$ws.Cells.Item(31, 14).Value = -795225

# Row 34
$ws.Cells.Item(34, 8).Value = 296675.16
$ws.Cells.Item(34, 9).Value = 38473.742
$ws.Cells.Item(34, 10).Value = 794635
$ws.Cells.Item(34, 11).Value = 38473.742
$ws.Cells.Item(34, 12).Value = 794635
$ws.Cells.Item(34, 13).Value = -38271.742
$ws.Cells.Item(34, 14).Value = -795039

# Row 86
$ws.Cells.Item(86, 8).Value = 3067.5557
$ws.Cells.Item(86, 9).Value = 2400
$ws.Cells.Item(86, 11).Value = 2400
$ws.Cells.Item(86, 13).Value = -1277

# Row 89
$ws.Cells.Item(89, 8).Value = 3067.5557
$ws.Cells.Item(89, 9).Value = 2400
$ws.Cells.Item(89, 11).Value = 12000
$ws.Cells.Item(89, 13).Value = -6384

# Row 107
$ws.Cells.Item(107, 8).Value = 411.44
$ws.Cells.Item(107, 9).Value = 367.54285
$ws.Cells.Item(107, 10).Value = 513.86664
$ws.Cells.Item(107, 11).Value = 367.54285
$ws.Cells.Item(107, 12).Value = 513.86664
$ws.Cells.Item(107, 13).Value = 1552.45715
$ws.Cells.Item(107, 14).Value = -4353.86664


$ws = $wb.Worksheets.Item("CUL")
# Row 13
$ws.Cells.Item(13, 8).Value = 2381.6667
$ws.Cells.Item(13, 9).Value = 145.5
$ws.Cells.Item(13, 10).Value = 3499.75
$ws.Cells.Item(13, 11).Value = 436.5
$ws.Cells.Item(13, 12).Value = 10499.25
$ws.Cells.Item(13, 13).Value = -268.5
$ws.Cells.Item(13, 14).Value = -10835.25

# Row 82
$ws.Cells.Item(82, 8).Value = 3700
$ws.Cells.Item(82, 9).Value = 0
$ws.Cells.Item(82, 10).Value = 3700
$ws.Cells.Item(82, 11).Value = 0
$ws.Cells.Item(82, 12).Value = 11100
$ws.Cells.Item(82, 13).Value = ""
$ws.Cells.Item(82, 14).Value = -11912

# Row 85
$ws.Cells.Item(85, 8).Value = 3700
$ws.Cells.Item(85, 9).Value = 0
$ws.Cells.Item(85, 10).Value = 3700
$ws.Cells.Item(85, 11).Value = 0
$ws.Cells.Item(85, 12).Value = 11100
$ws.Cells.Item(85, 13).Value = ""
$ws.Cells.Item(85, 14).Value = -13908

# Row 88
$ws.Cells.Item(88, 8).Value = 2966.6667
$ws.Cells.Item(88, 10).Value = 2966.6667
$ws.Cells.Item(88, 12).Value = 8900.000100000001
$ws.Cells.Item(88, 14).Value = -9756.000100000001

# Row 91
$ws.Cells.Item(91, 8).Value = 2966.6667
$ws.Cells.Item(91, 10).Value = 2966.6667
$ws.Cells.Item(91, 12).Value = 8900.000100000001
$ws.Cells.Item(91, 14).Value = -11864.0001

# Row 131
$ws.Cells.Item(131, 8).Value = 827.4706
$ws.Cells.Item(131, 10).Value = 1041.7916
$ws.Cells.Item(131, 12).Value = 3125.3748
$ws.Cells.Item(131, 14).Value = -13205.3748

# Row 133
$ws.Cells.Item(133, 8).Value = 8500
$ws.Cells.Item(133, 9).Value = 11000
$ws.Cells.Item(133, 10).Value = 7944.4443
$ws.Cells.Item(133, 11).Value = 33000
$ws.Cells.Item(133, 12).Value = 23833.3329
$ws.Cells.Item(133, 13).Value = -27940
$ws.Cells.Item(133, 14).Value = -33953.3329

# Row 137
$ws.Cells.Item(137, 8).Value = 1752.762
$ws.Cells.Item(137, 10).Value = 2414.8333
$ws.Cells.Item(137, 12).Value = 7244.499899999999
$ws.Cells.Item(137, 14).Value = -17444.4999

# Row 141
$ws.Cells.Item(141, 8).Value = 8866.105
$ws.Cells.Item(141, 9).Value = 2888
$ws.Cells.Item(141, 10).Value = 19114.285
$ws.Cells.Item(141, 11).Value = 8664
$ws.Cells.Item(141, 12).Value = 57342.855
$ws.Cells.Item(141, 13).Value = -3484
$ws.Cells.Item(141, 14).Value = -67702.855


$ws = $wb.Worksheets.Item("GSM")
# Row 38
$ws.Cells.Item(38, 8).Value = 0
$ws.Cells.Item(38, 10).Value = 0
$ws.Cells.Item(38, 12).Value = 0
$ws.Cells.Item(38, 14).Value = ""


$ws = $wb.Worksheets.Item("LTW")
# Row 60
$ws.Cells.Item(60, 8).Value = 11000
$ws.Cells.Item(60, 10).Value = 11000
$ws.Cells.Item(60, 12).Value = 11000
$ws.Cells.Item(60, 14).Value = -12018

# Row 61
$ws.Cells.Item(61, 8).Value = 2078.9524
$ws.Cells.Item(61, 9).Value = 2188.3076
$ws.Cells.Item(61, 10).Value = 1901.25
$ws.Cells.Item(61, 11).Value = 2188.3076
$ws.Cells.Item(61, 12).Value = 1901.25
$ws.Cells.Item(61, 13).Value = -1986.3076
$ws.Cells.Item(61, 14).Value = -2305.25

# Row 82
$ws.Cells.Item(82, 8).Value = 1139
$ws.Cells.Item(82, 9).Value = 1139
$ws.Cells.Item(82, 10).Value = 0
$ws.Cells.Item(82, 11).Value = 1139
$ws.Cells.Item(82, 12).Value = 0
$ws.Cells.Item(82, 13).Value = -778
$ws.Cells.Item(82, 14).Value = ""

# Row 85
$ws.Cells.Item(85, 8).Value = 1139
$ws.Cells.Item(85, 9).Value = 1139
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 11).Value = 1139
$ws.Cells.Item(85, 12).Value = 0
$ws.Cells.Item(85, 13).Value = 109
$ws.Cells.Item(85, 14).Value = ""

# Row 93
$ws.Cells.Item(93, 8).Value = 1263.9584
$ws.Cells.Item(93, 9).Value = 1216.75
$ws.Cells.Item(93, 10).Value = 1500
$ws.Cells.Item(93, 11).Value = 1216.75
$ws.Cells.Item(93, 12).Value = 1500
$ws.Cells.Item(93, 13).Value = 31.25
$ws.Cells.Item(93, 14).Value = -3996

# Row 100
$ws.Cells.Item(100, 8).Value = 1804.5625
$ws.Cells.Item(100, 9).Value = 1654.7778
$ws.Cells.Item(100, 11).Value = 1654.7778
$ws.Cells.Item(100, 13).Value = -1113.7778

# Row 113
$ws.Cells.Item(113, 8).Value = 2078.9524
$ws.Cells.Item(113, 9).Value = 2188.3076
$ws.Cells.Item(113, 10).Value = 1901.25
$ws.Cells.Item(113, 11).Value = 2188.3076
$ws.Cells.Item(113, 12).Value = 1901.25
$ws.Cells.Item(113, 13).Value = -18.30760000000009
$ws.Cells.Item(113, 14).Value = -6241.25


$ws = $wb.Worksheets.Item("WVR")
# Row 94
$ws.Cells.Item(94, 8).Value = 0
$ws.Cells.Item(94, 10).Value = 0
$ws.Cells.Item(94, 12).Value = 0
$ws.Cells.Item(94, 14).Value = ""

# Row 107
$ws.Cells.Item(107, 8).Value = 450.93332
$ws.Cells.Item(107, 9).Value = 384.04544
$ws.Cells.Item(107, 10).Value = 634.875
$ws.Cells.Item(107, 11).Value = 1152.13632
$ws.Cells.Item(107, 12).Value = 1904.625
$ws.Cells.Item(107, 13).Value = 767.8636799999999
$ws.Cells.Item(107, 14).Value = -5744.625

